$wb = $excel.ActiveWorkbook

# --- Sheet "Version": update Start Time, End Time, Run Time strings ---
$wsVersion = $wb.Worksheets.Item("Version")
$wsVersion.Range("C10").Value = "21:36:18 03-21-2017 Tuesday Pacific Daylight Time"
$wsVersion.Range("C11").Value = "21:36:53 03-21-2017 Tuesday Pacific Daylight Time"
$wsVersion.Range("C12").Value = "0:0:34 0 days"

# --- Sheet "3_3": update Deep Sleep (row3) and BT Idle (row4) stats ---
$ws33 = $wb.Worksheets.Item("3_3")
$ws33.Range("C3").Value = 0.581
$ws33.Range("D3").Value = 0.606
$ws33.Range("E3").Value = 0.555
$ws33.Range("H3").Value = "0.592485,0.565474,0.594074,0.566968,0.589244,0.593415,0.593327,0.593327,0.591846,0.566221,0.596024,0.572444,0.562695,0.590175,0.56919,0.59063,0.593409,0.563802,0.591187,0.56046,0.587655,0.5951,0.594726,0.563619,0.58794,0.594808,0.59546,0.589896,0.592131,0.554991,0.586635,0.558504,0.595739,0.563626,0.594448,0.561024,0.592953,0.569183,0.588136,0.567695,0.600379,0.564448,0.592396,0.591839,0.56724,0.592111,0.555256,0.593517,0.557124,0.604557,0.559543,0.592491,0.555643,0.596764,0.561112,0.598334,0.594346,0.565942,0.584509,0.562878,0.598062,0.559529,0.596846,0.569658,0.597315,0.556744,0.597593,0.569747,0.586173,0.558232,0.595182,0.600284,0.562416,0.592858,0.564828,0.594257,0.561207,0.592219,0.563986,0.597885,0.56046,0.588598,0.563619,0.593612,0.560745,0.593327,0.588497,0.605861,0.594991,0.566955,0.594529,0.562885,0.584135,0.564088,0.598252,0.566221,0.588028,0.563259,0.600113,0.565283"

$ws33.Range("C4").Value = 9.798999999999999
$ws33.Range("D4").Value = 9.818
$ws33.Range("E4").Value = 9.779
$ws33.Range("F4").Value = 0.008
$ws33.Range("H4").Value = "9.793752,9.802672,9.810648,9.800165,9.795709,9.791246,9.797183,9.807224,9.795607,9.782604,9.792448,9.810186,9.807584,9.797835,9.784934,9.795614,9.810464,9.800247,9.788739,9.797183,9.800994,9.809072,9.79988,9.786334,9.797842,9.805267,9.80509,9.789581,9.785865,9.796911,9.8113,9.801646,9.78558,9.786599,9.807387,9.809262,9.802482,9.779269,9.801735,9.810091,9.80916,9.790974,9.786612,9.79681,9.803318,9.796531,9.786137,9.79329,9.803026,9.807686,9.794962,9.788087,9.794676,9.811008,9.812034,9.789941,9.791429,9.805179,9.804527,9.808691,9.786707,9.795702,9.802102,9.8177,9.792353,9.788657,9.790125,9.811762,9.809819,9.793385,9.786796,9.79524,9.80935,9.803419,9.790505,9.79776,9.808882,9.808229,9.799703,9.794119,9.805444,9.805356,9.798222,9.791803,9.795512,9.802298,9.816687,9.79348,9.78846,9.794309,9.804602,9.803039,9.797747,9.791708,9.797835,9.810926,9.80397,9.793847,9.794139,9.804432"

# --- Sheet "1_8": update Deep Sleep (row3) and BT Idle (row4) stats ---
$ws18 = $wb.Worksheets.Item("1_8")
$ws18.Range("C3").Value = -0.002
$ws18.Range("D3").Value = 0.006
$ws18.Range("E3").Value = -0.013
$ws18.Range("F3").Value = 0.003
$ws18.Range("H3").Value = "0.000613,-0.002742,-0.004129,-0.002187,0.001704,-0.001385,-0.005264,-0.001897,-0.00691,-0.001367,-0.003007,-0.002168,-0.00247,-0.00022,0.000311,0.005577,-0.002458,-0.001077,-0.005239,-0.005794,-0.013286,-0.006343,-0.004407,0.000298,-0.006608,-0.004672,-0.004678,0.000866,4.5e-05,-0.003013,-0.002452,-0.003013,-0.009395,-0.001613,0.00393,-0.002433,-0.002168,-0.00136,-0.002729,0.003104,0.000317,-0.004974,-0.005825,-0.002174,-0.003574,-0.000805,-0.002458,0.001686,-0.004123,-0.006626,-0.000793,-0.002187,0.003369,-0.002729,3.9e-05,-0.001626,0.003918,-0.000232,-0.003562,-0.000818,0.001162,-0.005788,-0.006078,-0.000787,-0.007459,-0.005794,-0.001891,-0.001891,0.001704,-0.00133,-0.003556,-0.003833,-0.003562,-0.003858,-0.004407,0.003085,-0.001613,3.3e-05,-0.00276,-0.010795,-0.001891,-0.005214,0.001439,-0.003303,0.001131,0.001433,-0.002181,-0.006639,-0.001626,0.003344,0.000576,-0.003852,-0.007169,-0.001613,0.006144,-0.002174,-0.000238,-0.005245,-0.005529,-0.006608"

$ws18.Range("C4").Value = -0.002
$ws18.Range("D4").Value = 0.005
$ws18.Range("E4").Value = -0.008999999999999999
$ws18.Range("F4").Value = 0.003
$ws18.Range("H4").Value = "-0.000812,0.00393,-0.00498,-0.004702,-0.007181,0.000323,0.001421,0.000847,-0.007742,-0.000787,-0.00247,-0.0044,-0.003568,-0.002711,-0.002476,0.004756,-0.00638,0.004775,-0.001638,-0.001891,-0.001922,0.000329,3.3e-05,0.003659,0.004183,-0.000509,-0.007163,-0.001077,-0.006337,-0.003019,0.000853,-0.005523,0.000872,-0.00329,-0.000793,-0.001058,0.004485,-0.0044,-0.008291,-0.00329,-0.004117,-0.004659,-0.002168,-0.004955,-0.003593,-0.003821,0.002259,-0.000244,-0.003574,0.004195,-0.001607,-0.004918,-0.003833,0.003351,-0.006608,-0.002174,-0.001058,-0.003839,-0.004918,-0.004943,-0.005782,0.002259,-0.002729,-0.000805,-0.00609,2.1e-05,0.001445,0.002253,-0.000516,-0.000497,0.000884,-0.001077,-0.005782,-0.009142,-0.005486,-0.002476,-0.004129,0.001988,0.000563,-0.000793,-0.006916,0.004491,-0.003284,-0.004135,-0.000534,0.002247,-0.005775,-0.006361,0.001137,0.001975,-0.002458,0.005324,-0.000818,-0.003858,-0.004653,-0.004672,-0.00136,0.000884,-0.001077,0.003085"
